$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 333
$ws.Cells.Item(3, 6).Value = 1152
$ws.Cells.Item(4, 6).Value = 607
$ws.Cells.Item(7, 6).Value = 162
$ws.Cells.Item(8, 6).Value = 661
$ws.Cells.Item(9, 6).Value = 1858
$ws.Cells.Item(10, 6).Value = 58
$ws.Cells.Item(11, 6).Value = 441
$ws.Cells.Item(12, 6).Value = 65
$ws.Cells.Item(13, 6).Value = 81
$ws.Cells.Item(14, 6).Value = 698
$ws.Cells.Item(15, 6).Value = 465
$ws.Cells.Item(17, 6).Value = 823
$ws.Cells.Item(18, 6).Value = 80549
$ws.Cells.Item(19, 6).Value = 80549
$ws.Cells.Item(22, 6).Value = 33988
$ws.Cells.Item(23, 6).Value = 33988
$ws.Cells.Item(24, 6).Value = 559
$ws.Cells.Item(26, 6).Value = 31
$ws.Cells.Item(27, 6).Value = 65
$ws.Cells.Item(28, 6).Value = 59
$ws.Cells.Item(29, 6).Value = 1021
$ws.Cells.Item(30, 6).Value = 323
$ws.Cells.Item(31, 6).Value = 165
$ws.Cells.Item(32, 6).Value = 661
$ws.Cells.Item(33, 6).Value = 3053
$ws.Cells.Item(34, 6).Value = 3053
$ws.Cells.Item(35, 6).Value = 1238
$ws.Cells.Item(36, 6).Value = 5522
$ws.Cells.Item(37, 6).Value = 817
$ws.Cells.Item(38, 6).Value = 470
$ws.Cells.Item(41, 6).Value = 6
$ws.Cells.Item(42, 6).Value = 457
$ws.Cells.Item(46, 6).Value = 60

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(10, 6).Value = 1995
$ws.Cells.Item(11, 6).Value = 36
$ws.Cells.Item(13, 6).Value = 89
$ws.Cells.Item(14, 6).Value = 421
$ws.Cells.Item(24, 6).Value = 27
$ws.Cells.Item(27, 6).Value = 83
$ws.Cells.Item(31, 6).Value = 1676
$ws.Cells.Item(32, 6).Value = 500
$ws.Cells.Item(43, 6).Value = 74
$ws.Cells.Item(44, 6).Value = 832
$ws.Cells.Item(45, 6).Value = 239

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(5, 6).Value = 588
$ws.Cells.Item(6, 6).Value = 618
$ws.Cells.Item(7, 6).Value = 184

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 333
$ws.Cells.Item(4, 6).Value = 1152
$ws.Cells.Item(5, 6).Value = 618
$ws.Cells.Item(6, 6).Value = 618
$ws.Cells.Item(9, 6).Value = 607
$ws.Cells.Item(11, 6).Value = 162
$ws.Cells.Item(12, 6).Value = 661
$ws.Cells.Item(13, 6).Value = 184
$ws.Cells.Item(14, 6).Value = 1858
$ws.Cells.Item(15, 6).Value = 36
$ws.Cells.Item(16, 6).Value = 58
$ws.Cells.Item(17, 6).Value = 441
$ws.Cells.Item(18, 6).Value = 65
$ws.Cells.Item(19, 6).Value = 81
$ws.Cells.Item(20, 6).Value = 824
$ws.Cells.Item(23, 6).Value = 80549
$ws.Cells.Item(25, 6).Value = 33988
$ws.Cells.Item(26, 6).Value = 559
$ws.Cells.Item(28, 6).Value = 31
$ws.Cells.Item(31, 6).Value = 59
$ws.Cells.Item(34, 6).Value = 323
$ws.Cells.Item(35, 6).Value = 165
$ws.Cells.Item(36, 6).Value = 27
$ws.Cells.Item(37, 6).Value = 3053
$ws.Cells.Item(38, 6).Value = 1238
$ws.Cells.Item(39, 6).Value = 5522
$ws.Cells.Item(40, 6).Value = 83
$ws.Cells.Item(41, 6).Value = 817
$ws.Cells.Item(43, 6).Value = 1676
$ws.Cells.Item(46, 6).Value = 6
$ws.Cells.Item(47, 6).Value = 457
$ws.Cells.Item(50, 6).Value = 74
$ws.Cells.Item(52, 6).Value = 239
$ws.Cells.Item(54, 6).Value = 60
